$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("H70").Value = 3288.3333
$ws.Range("I70").Value = 2687
$ws.Range("J70").Value = 3769.4
$ws.Range("K70").Value = 8061
$ws.Range("L70").Value = 11308.2
$ws.Range("M70").Value = -7791
$ws.Range("N70").Value = -11848.2
$ws.Range("H73").Value = 3288.3333
$ws.Range("I73").Value = 2687
$ws.Range("J73").Value = 3769.4
$ws.Range("K73").Value = 8061
$ws.Range("L73").Value = 11308.2
$ws.Range("M73").Value = -7125
$ws.Range("N73").Value = -13180.2
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H107").Value = 2762.4
$ws.Range("I107").Value = 1255.3334
$ws.Range("J107").Value = 5023
$ws.Range("K107").Value = 1255.3334
$ws.Range("L107").Value = 5023
$ws.Range("M107").Value = 664.6666
$ws.Range("N107").Value = -8863
$ws.Range("H137").Value = 1595.2413
$ws.Range("J137").Value = 2242
$ws.Range("L137").Value = 6726
$ws.Range("N137").Value = -11826
$ws.Range("H138").Value = 3811.4285
$ws.Range("J138").Value = 3811.4285
$ws.Range("L138").Value = 11434.2855
$ws.Range("N138").Value = -21714.2855
$ws.Range("H141").Value = 5065
$ws.Range("I141").Value = 4712.7144
$ws.Range("K141").Value = 14138.1432
$ws.Range("M141").Value = -8958.143199999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5176.2095
$ws.Range("I32").Value = 4062.8538
$ws.Range("K32").Value = 4062.8538
$ws.Range("M32").Value = -3775.8538
$ws.Range("H37").Value = 20833.334
$ws.Range("J37").Value = 25000
$ws.Range("L37").Value = 25000
$ws.Range("N37").Value = -25546
$ws.Range("H44").Value = 35000
$ws.Range("J44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("N44").Value = -35976
$ws.Range("H55").Value = 31666.666
$ws.Range("J55").Value = 31666.666
$ws.Range("L55").Value = 31666.666
$ws.Range("N55").Value = -32296.666
$ws.Range("H61").Value = 2229.6155
$ws.Range("I61").Value = 1387.2222
$ws.Range("K61").Value = 1387.2222
$ws.Range("M61").Value = -1175.2222
$ws.Range("H74").Value = 4171.875
$ws.Range("I74").Value = 4171.875
$ws.Range("K74").Value = 4171.875
$ws.Range("M74").Value = -3297.875
$ws.Range("H77").Value = 4171.875
$ws.Range("I77").Value = 4171.875
$ws.Range("K77").Value = 20859.375
$ws.Range("M77").Value = -16491.375
$ws.Range("H80").Value = 32498.75
$ws.Range("I80").Value = 9995
$ws.Range("K80").Value = 9995
$ws.Range("M80").Value = -8997
$ws.Range("H83").Value = 32498.75
$ws.Range("I83").Value = 9995
$ws.Range("K83").Value = 29985
$ws.Range("M83").Value = -24993
$ws.Range("H102").Value = 770
$ws.Range("I102").Value = 770
$ws.Range("K102").Value = 770
$ws.Range("M102").Value = 852
$ws.Range("H110").Value = 2891.0625
$ws.Range("I110").Value = 782.25
$ws.Range("K110").Value = 782.25
$ws.Range("M110").Value = 1262.75
$ws.Range("H136").Value = 2229.6155
$ws.Range("I136").Value = 1387.2222
$ws.Range("K136").Value = 4161.6666
$ws.Range("M136").Value = -1611.6666
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H105").Value = 3124.75
$ws.Range("I105").Value = 2999.5
$ws.Range("K105").Value = 2999.5
$ws.Range("M105").Value = -1252.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 2295
$ws.Range("I45").Value = 2295
$ws.Range("K45").Value = 2295
$ws.Range("M45").Value = -1702
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1265
$ws.Range("J75").Value = 2500
$ws.Range("L75").Value = 7500
$ws.Range("N75").Value = -9496
$ws.Range("H78").Value = 1265
$ws.Range("J78").Value = 2500
$ws.Range("L78").Value = 22500
$ws.Range("N78").Value = -32484
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 830.06665
$ws.Range("I97").Value = 832.2143
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 832.2143
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -336.2143
$ws.Range("N97").Value = -1792
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2299.7778
$ws.Range("I46").Value = 800
$ws.Range("K46").Value = 800
$ws.Range("M46").Value = -612
$ws.Range("H56").Value = 42500
$ws.Range("I56").Value = 60000
$ws.Range("J56").Value = 25000
$ws.Range("K56").Value = 60000
$ws.Range("L56").Value = 25000
$ws.Range("M56").Value = -59309
$ws.Range("N56").Value = -26382
$ws.Range("H76").Value = 69999
$ws.Range("J76").Value = 69999
$ws.Range("L76").Value = 69999
$ws.Range("N76").Value = -70675
$ws.Range("H79").Value = 69999
$ws.Range("J79").Value = 69999
$ws.Range("L79").Value = 69999
$ws.Range("N79").Value = -72339
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 46029
$ws.Range("J37").Value = 46029
$ws.Range("L37").Value = 46029
$ws.Range("N37").Value = -46435
$ws.Range("H81").Value = 569.6
$ws.Range("I81").Value = 462
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 924
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = 137
$ws.Range("N81").Value = -4122
$ws.Range("H84").Value = 569.6
$ws.Range("I84").Value = 462
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 4620
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = 684
$ws.Range("N84").Value = -20608
$ws.Range("H96").Value = 1020
$ws.Range("I96").Value = 775
$ws.Range("K96").Value = 775
$ws.Range("M96").Value = 598
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
